$p = $ppt.ActivePresentation

# Add a new slide at the end, using the "Title and Content" layout
# (the same layout used by the existing content slides, slideLayout2.xml,
# which is the 2nd layout in the slide master's layout list -> ppLayoutText).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

# Title placeholder
$s.Shapes.Item(1).TextFrame.TextRange.Text = "References"

# Body / content placeholder
$tf = $s.Shapes.Item(2).TextFrame
$url = "https://julie-jiang.github.io/image-segmentation/"
$mid = " (useful explanation but same algorithms as those in prof. Brower" + [char]8217 + "s "
$tail = "example presentation)"

$tf.TextRange.Text = $url
$r2 = $tf.TextRange.InsertAfter($mid)
$r3 = $r2.InsertAfter($tail)

$linkRange = $tf.TextRange.Characters(1, $url.Length)
$linkRange.ActionSettings(1).Hyperlink.Address = $url
